$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "25.940.47"
$ws.Cells.Item(2, 5).Value = "  +0.67%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.634.33"
$ws.Cells.Item(3, 5).Value = "  +0.25%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.26%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'214.77"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.07%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'0.505"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +0.85%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.21%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.09%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  -0.08%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  +0.53%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.0792"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.00%  "

# Row 12
$ws.Cells.Item(12, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(12, 4).Value = "1.860.13"
$ws.Cells.Item(12, 5).Value = "  +0.19%  "

# Row 13
$ws.Cells.Item(13, 2).Value = "Polkadot"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(13, 4).Value = "'4.24"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -0.54%  "

# Row 14
$ws.Cells.Item(14, 2).Value = "WrappedEther"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(14, 4).Value = "1.627.05"
$ws.Cells.Item(14, 5).Value = "  -2.27%  "

# Row 15
$ws.Cells.Item(15, 5).Value = "  -1.61%  "

# Row 16
$ws.Cells.Item(16, 2).Value = "ShibaInu"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(16, 4).Value = "0.0₃0756"
$ws.Cells.Item(16, 5).Value = "  -0.69%  "

# Row 17
$ws.Cells.Item(17, 2).Value = "Litecoin"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(17, 4).Value = "'62.89"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -0.05%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "25.925.89"
$ws.Cells.Item(18, 5).Value = "  +0.65%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  +0.30%  "

# Row 20
$ws.Cells.Item(20, 2).Value = "BitcoinCash"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(20, 4).Value = "'193.04"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +0.94%  "

# Row 21
$ws.Cells.Item(21, 2).Value = "Uniswap"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(21, 4).Value = "'4.38"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -1.38%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'9.98"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.58%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  -0.09%  "

# Row 25
$ws.Cells.Item(25, 2).Value = "BinanceUSD"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(25, 4).Value = "'1.00"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +0.09%  "

# Row 26
$ws.Cells.Item(26, 2).Value = "Monero"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(26, 4).Value = "'142.74"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +0.12%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  +1.75%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +0.40%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'15.47"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -0.07%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  +0.07%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'0.0498"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +0.63%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'3.30"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -0.66%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'3.23"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -0.42%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  -0.68%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  +2.12%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -0.56%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "1.139.10"
$ws.Cells.Item(37, 5).Value = "  +0.01%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'0.550"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +1.26%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  -1.28%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  +0.57%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +0.29%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "TrustWalletToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(42, 4).Value = "'0.804"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -0.05%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "FraxShare"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(43, 4).Value = "'5.47"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -1.25%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "'99.25"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -1.55%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "1.769.46"
$ws.Cells.Item(45, 5).Value = "  +0.18%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "0.0₆0111"
$ws.Cells.Item(46, 5).Value = "  -0.10%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "'56.28"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +1.90%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "'0.0525"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +2.98%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  +0.76%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  -0.68%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "'7.64"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +2.33%  "
